# Update "想去人数" (number of people interested) counts on the
# "展览" and "全部类型" sheets to reflect the latest generated output.
#
# Changes:
#   F2:  138   -> 139
#   F7:  11706 -> 11708
#   F12: 1099  -> 1100
#   F15: 13254 -> 13255
# applied identically on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    "F2"  = 139
    "F7"  = 11708
    "F12" = 1100
    "F15" = 13255
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
